$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions) - update column F ("想去人数" / interest counts)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 3250
$ws.Cells.Item(4, 6).Value = 2009
$ws.Cells.Item(7, 6).Value = 3093
$ws.Cells.Item(8, 6).Value = 618
$ws.Cells.Item(10, 6).Value = 40
$ws.Cells.Item(12, 6).Value = 155
$ws.Cells.Item(15, 6).Value = 10197
$ws.Cells.Item(17, 6).Value = 238
$ws.Cells.Item(18, 6).Value = 13
$ws.Cells.Item(19, 6).Value = 38
$ws.Cells.Item(20, 6).Value = 8099
$ws.Cells.Item(21, 6).Value = 12705
$ws.Cells.Item(24, 6).Value = 27
$ws.Cells.Item(25, 6).Value = 274
$ws.Cells.Item(27, 6).Value = 601
$ws.Cells.Item(28, 6).Value = 8
$ws.Cells.Item(29, 6).Value = 424
$ws.Cells.Item(32, 6).Value = 241
$ws.Cells.Item(33, 6).Value = 7995
$ws.Cells.Item(34, 6).Value = 1582
$ws.Cells.Item(35, 6).Value = 221
$ws.Cells.Item(37, 6).Value = 86
$ws.Cells.Item(38, 6).Value = 4628
$ws.Cells.Item(39, 6).Value = 1456
$ws.Cells.Item(41, 6).Value = 385
$ws.Cells.Item(42, 6).Value = 85
$ws.Cells.Item(43, 6).Value = 642

# Sheet 2: 演出 (Performances) - update column F
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(12, 6).Value = 29

# Sheet 3: 本地生活 (Local Life) - update column F
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 647
$ws.Cells.Item(5, 6).Value = 21

# Sheet 4: 全部类型 (All Types, combined view) - update column F
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 647
$ws.Cells.Item(4, 6).Value = 3250
$ws.Cells.Item(6, 6).Value = 2009
$ws.Cells.Item(9, 6).Value = 21
$ws.Cells.Item(10, 6).Value = 3093
$ws.Cells.Item(12, 6).Value = 618
$ws.Cells.Item(13, 6).Value = 40
$ws.Cells.Item(15, 6).Value = 155
$ws.Cells.Item(18, 6).Value = 10197
$ws.Cells.Item(19, 6).Value = 238
$ws.Cells.Item(20, 6).Value = 13
$ws.Cells.Item(21, 6).Value = 38
$ws.Cells.Item(22, 6).Value = 8099
$ws.Cells.Item(23, 6).Value = 12705
$ws.Cells.Item(25, 6).Value = 27
$ws.Cells.Item(26, 6).Value = 274
$ws.Cells.Item(28, 6).Value = 601
$ws.Cells.Item(30, 6).Value = 8
$ws.Cells.Item(33, 6).Value = 29
$ws.Cells.Item(35, 6).Value = 241
$ws.Cells.Item(36, 6).Value = 7995
$ws.Cells.Item(37, 6).Value = 221
$ws.Cells.Item(39, 6).Value = 86
$ws.Cells.Item(40, 6).Value = 4628
$ws.Cells.Item(47, 6).Value = 642
